$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "240.54" or "1.00"
# are not auto-converted to numbers by Excel, matching the original inline-string formatting.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '36.360.10'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '1.930.42'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '240.54'
$ws.Range('E5').Value = '  -2.11%  '
$ws.Range('D6').Value = '0.604'
$ws.Range('E6').Value = '  -2.93%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '55.76'
$ws.Range('E8').Value = '  -5.55%  '
$ws.Range('D9').Value = '0.355'
$ws.Range('E9').Value = '  -5.61%  '
$ws.Range('D10').Value = '0.0832'
$ws.Range('E10').Value = '  +0.86%  '
$ws.Range('E11').Value = '  -0.96%  '
$ws.Range('D12').Value = '2.214.60'
$ws.Range('E12').Value = '  -2.70%  '
$ws.Range('D13').Value = '0.794'
$ws.Range('E13').Value = '  -8.36%  '
$ws.Range('D14').Value = '13.22'
$ws.Range('E14').Value = '  -5.76%  '
$ws.Range('D15').Value = '20.64'
$ws.Range('E15').Value = '  -12.22%  '
$ws.Range('D16').Value = '5.08'
$ws.Range('E16').Value = '  -7.12%  '
$ws.Range('D17').Value = '1.936.52'
$ws.Range('E17').Value = '  -2.33%  '
$ws.Range('D18').Value = '36.273.72'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').Value = '68.52'
$ws.Range('E19').Value = '  -2.80%  '
$ws.Range('D20').Value = '0.0₃0854'
$ws.Range('E20').Value = '  -2.62%  '
$ws.Range('D21').Value = '225.66'
$ws.Range('E21').Value = '  -3.69%  '
$ws.Range('D22').Value = '4.90'
$ws.Range('E22').Value = '  -7.78%  '
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').Value = '2.31'
$ws.Range('E24').Value = '  -10.40%  '
$ws.Range('E25').Value = '  -2.85%  '
$ws.Range('D26').Value = '9.07'
$ws.Range('E26').Value = '  -9.10%  '
$ws.Range('D27').Value = '160.16'
$ws.Range('E27').Value = '  -1.33%  '
$ws.Range('E28').Value = '  -2.67%  '
$ws.Range('D29').Value = '18.98'
$ws.Range('E29').Value = '  -4.68%  '
$ws.Range('E30').Value = '  -3.20%  '
$ws.Range('E31').Value = '  -7.93%  '
$ws.Range('D32').Value = '4.49'
$ws.Range('E32').Value = '  -8.53%  '
$ws.Range('D33').Value = '0.0617'
$ws.Range('E33').Value = '  -5.89%  '
$ws.Range('D34').Value = '4.11'
$ws.Range('E34').Value = '  -7.18%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').Value = '5.99'
$ws.Range('E36').Value = '  -4.04%  '
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('E38').Value = '  -6.76%  '
$ws.Range('E39').Value = '  -0.85%  '
$ws.Range('E40').Value = '  -0.76%  '
$ws.Range('E41').Value = '  -1.19%  '
$ws.Range('D42').Value = '0.0208'
$ws.Range('E42').Value = '  -2.78%  '
$ws.Range('E43').Value = '  -8.26%  '
$ws.Range('D44').Value = '15.38'
$ws.Range('E44').Value = '  -5.41%  '
$ws.Range('D45').Value = '1.325.00'
$ws.Range('E45').Value = '  -3.03%  '
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -8.03%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '7.02'
$ws.Range('E47').Value = '  -6.35%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '84.36'
$ws.Range('E48').Value = '  -8.88%  '
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('D50').Value = '2.105.91'
$ws.Range('E50').Value = '  -2.70%  '
$ws.Range('D51').Value = '43.09'

# Restore default (General) style so the cells keep the workbook's original formatting/style index
$priceRange.ClearFormats()
